{"js": "// Remove the inline picture from the document body, leaving an empty\n// paragraph where it used to live (the rest of the document is unchanged).\nconst body = context.document.body;\nconst pictures = body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < pictures.items.length; i++) {\n  pictures.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the inline picture from the document, leaving an empty paragraph\n# where it used to live (the rest of the document is unchanged).\n$d = $word.ActiveDocument\n\nfor ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {\n    $d.InlineShapes.Item($i).Delete()\n}\n"}
